# ObjectLocators.xlsx -- add "PIMPage" sheet (Add Employee / AutoIt locators)
# and a new "pimDashboardLocator" row on DashboardPage.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new worksheet "PIMPage" right after "DashboardPage" (it ends up
#    between DashboardPage and Locators, matching the workbook.xml sheet order).
# ---------------------------------------------------------------------------
$dashSheet = $wb.Worksheets.Item("DashboardPage")
$pimSheet = $wb.Worksheets.Add($null, $dashSheet)
$pimSheet.Name = "PIMPage"

# ---------------------------------------------------------------------------
# 2. Populate column-A labels for rows 1-9 first (matches shared-string order).
# ---------------------------------------------------------------------------
$pimSheet.Range("A1").Value = "addEmployeeNavigationMenuButton"
$pimSheet.Range("A2").Value = "firstNameTextBox"
$pimSheet.Range("A3").Value = "middleNameTextBox"
$pimSheet.Range("A4").Value = "lastNameTextBox"
$pimSheet.Range("A5").Value = "employeeIdTextBox"
$pimSheet.Range("A6").Value = "loginDetailsToggle"
$pimSheet.Range("A7").Value = "setUsernameTextBox"
$pimSheet.Range("A8").Value = "setPasswordTextBox"
$pimSheet.Range("A9").Value = "confirmPasswordTextBox"

# ---------------------------------------------------------------------------
# 3. Column-B type marker ("XPATH") for rows 1-9.
# ---------------------------------------------------------------------------
$pimSheet.Range("B1").Value = "XPATH"
$pimSheet.Range("B2").Value = "XPATH"
$pimSheet.Range("B3").Value = "XPATH"
$pimSheet.Range("B4").Value = "XPATH"
$pimSheet.Range("B5").Value = "XPATH"
$pimSheet.Range("B6").Value = "XPATH"
$pimSheet.Range("B7").Value = "XPATH"
$pimSheet.Range("B8").Value = "XPATH"
$pimSheet.Range("B9").Value = "XPATH"

# ---------------------------------------------------------------------------
# 4. Column-C locator strings for rows 1-5, 7, 8, 9 (row 6 filled in later).
# ---------------------------------------------------------------------------
$pimSheet.Range("C1").Value = "//a[@class='oxd-topbar-body-nav-tab-item' and text()='Add Employee']"
$pimSheet.Range("C2").Value = "//input[@name='firstName']"
$pimSheet.Range("C3").Value = "//input[@name='middleName']"
$pimSheet.Range("C4").Value = "//input[@name='lastName']"
$pimSheet.Range("C5").Value = "(//input[@class='oxd-input oxd-input--active'])[2]"
$pimSheet.Range("C7").Value = "(//div[@class='oxd-input-group oxd-input-field-bottom-space']/div/input[@class='oxd-input oxd-input--active'])[2]"
$pimSheet.Range("C8").Value = "(//input[@type='password'])[1]"
$pimSheet.Range("C9").Value = "(//input[@type='password'])[2]"

# ---------------------------------------------------------------------------
# 5. DashboardPage gains a "pimDashboardLocator" row pointing at the PIM menu.
# ---------------------------------------------------------------------------
$dashSheet.Range("A2").Value = "pimDashboardLocator"
$dashSheet.Range("B2").Value = "XPATH"
$dashSheet.Range("C2").Value = "//span[@class='oxd-text oxd-text--span oxd-main-menu-item--name' and text()='PIM']"

# ---------------------------------------------------------------------------
# 6. Back to PIMPage: fill in row 6's locator (added after the dashboard edit).
# ---------------------------------------------------------------------------
$pimSheet.Range("C6").Value = "//div[@class='oxd-switch-wrapper']/label"

# ---------------------------------------------------------------------------
# 7. Row 10: upload-profile-picture button, added last.
# ---------------------------------------------------------------------------
$pimSheet.Range("A10").Value = "uploadProfilePictureButton"
$pimSheet.Range("B10").Value = "XPATH"
$pimSheet.Range("C10").Value = "//button[@class='oxd-icon-button employee-image-action']"

# ---------------------------------------------------------------------------
# 8. Column widths (tuned to the engine's pixel-grid rounding so the saved
#    OOXML <col width> lands as close as possible to the captured values).
# ---------------------------------------------------------------------------
$dashSheet.Columns.Item(1).ColumnWidth = 19.65
$pimSheet.Columns.Item(1).ColumnWidth = 37.33
$pimSheet.Columns.Item(3).ColumnWidth = 61.1

# ---------------------------------------------------------------------------
# 9. Selections / active sheet.
# ---------------------------------------------------------------------------
$dashSheet.Range("C2").Select()
$pimSheet.Range("C10").Select()
$pimSheet.Activate()
